$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '28.798.60'
$ws.Range('E2').Value = '  +7.67%  '
$ws.Range('D3').Value = '1.813.11'
$ws.Range('E3').Value = '  +5.13%  '
$ws.Range('E4').Value = '  +0.24%  '
Set-TextValue 'D5' '250.00'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '0.4958'
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('D8').Value = '0.2788'
$ws.Range('E8').Value = '  +8.00%  '
$ws.Range('D9').Value = '0.06402'
$ws.Range('E9').Value = '  +3.33%  '
$ws.Range('D10').Value = '1.809.15'
$ws.Range('E10').Value = '  +4.87%  '
$ws.Range('D11').Value = '16.78'
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('D12').Value = '0.07147'
$ws.Range('E12').Value = '  +3.59%  '
$ws.Range('D13').Value = '0.6514'
$ws.Range('E13').Value = '  +7.22%  '
Set-TextValue 'D14' '83.70'
$ws.Range('E14').Value = '  +8.94%  '
$ws.Range('D15').Value = '4.703'
$ws.Range('E15').Value = '  +5.06%  '
$ws.Range('D16').Value = '28.784.79'
Set-TextValue 'D17' '1.000'
$ws.Range('E17').Value = '  +0.20%  '
Set-TextValue 'D18' '0.000007402'
$ws.Range('E18').Value = '  +3.54%  '
Set-TextValue 'D19' '1.000'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = '12.24'
$ws.Range('E20').Value = '  +7.00%  '
$ws.Range('D21').Value = '2.053.09'
$ws.Range('E21').Value = '  +5.28%  '
$ws.Range('D22').Value = '4.612'
$ws.Range('E22').Value = '  +4.04%  '
Set-TextValue 'D23' '8.900'
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('D24').Value = '5.355'
$ws.Range('E24').Value = '  +5.66%  '
$ws.Range('D25').Value = '143.32'
$ws.Range('E25').Value = '  +4.62%  '
$ws.Range('D26').Value = '130.63'
$ws.Range('E26').Value = '  +23.35%  '
$ws.Range('D27').Value = '16.34'
$ws.Range('E27').Value = '  +7.12%  '
$ws.Range('D28').Value = '1.891'
$ws.Range('E28').Value = '  +6.83%  '
$ws.Range('D29').Value = '1.406'
$ws.Range('E29').Value = '  +1.91%  '
$ws.Range('D30').Value = '4.169'
$ws.Range('E30').Value = '  +5.90%  '
Set-TextValue 'D31' '0.08370'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('D32').Value = '3.862'
$ws.Range('E32').Value = '  +4.52%  '
$ws.Range('D33').Value = '0.04957'
$ws.Range('E33').Value = '  +10.39%  '
$ws.Range('E34').Value = '  +8.15%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D35' '2.720'
$ws.Range('E35').Value = '  +4.71%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6786'
$ws.Range('E36').Value = '  +9.50%  '
$ws.Range('D37').Value = '2.775'
$ws.Range('E37').Value = '  +13.90%  '
Set-TextValue 'D38' '2.240'
$ws.Range('E38').Value = '  +10.32%  '
$ws.Range('D39').Value = '0.9617'
$ws.Range('E39').Value = '  +3.75%  '
$ws.Range('D40').Value = '6.068'
$ws.Range('E40').Value = '  +7.23%  '
$ws.Range('D41').Value = '0.01593'
$ws.Range('E41').Value = '  +6.48%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = '100.85'
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('D44').Value = '0.4104'
$ws.Range('D45').Value = '7.232'
$ws.Range('E45').Value = '  +5.49%  '
$ws.Range('D46').Value = '0.1225'
$ws.Range('E46').Value = '  +5.94%  '
$ws.Range('D47').Value = '0.05515'
$ws.Range('E47').Value = '  +2.27%  '
$ws.Range('D48').Value = '8.228'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('D49').Value = '31.72'
$ws.Range('E49').Value = '  +5.35%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.312'
$ws.Range('E50').Value = '  +6.77%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '0.3634'
$ws.Range('E51').Value = '  +8.11%  '
